$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix DTR entries: the 0.5 "sick leave" values recorded for 04-28-2015 and
#    05-08-2015 were excess/incorrect; clear them out of the SICK LEAVE
#    column (I) and record the correct "undertime" hours in column F.
# ---------------------------------------------------------------------------
$ws.Range("F8").Value = 2.5
$ws.Range("I8").Value = $null

$ws.Range("F18").Value = 1.5
$ws.Range("I18").Value = $null

# ---------------------------------------------------------------------------
# 2. Add a "Legends:" section title above the color-coded legend rows, using
#    the same bold/underlined/size-15 Arial look as the report's title rows.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$legendTitle = $ws.Range("E24:P24")
$legendTitle.Merge() | Out-Null
$legendTitle.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$legendTitle.Value = "Legends:"

# ---------------------------------------------------------------------------
# 3. Legend row 1: "LATES" swatch (teal) + description.
# ---------------------------------------------------------------------------
$swatch1 = $ws.Range("E25:E26")
$swatch1.Merge() | Out-Null
$swatch1.Interior.Color = 0xCCA329

$desc1 = $ws.Range("F25:P26")
$desc1.Merge() | Out-Null
$desc1.Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."

# ---------------------------------------------------------------------------
# 4. Legend row 2: half-day swatch (orange) + description.
# ---------------------------------------------------------------------------
$swatch2 = $ws.Range("E27:E28")
$swatch2.Merge() | Out-Null
$swatch2.Interior.Color = 0x66CCFF

$desc2 = $ws.Range("F27:P28")
$desc2.Merge() | Out-Null
$desc2.Value = "Employee is considered half-day because of his time-in or time-out."

# ---------------------------------------------------------------------------
# 5. Legend row 3: absent swatch (red) + description.
# ---------------------------------------------------------------------------
$swatch3 = $ws.Range("E29:E30")
$swatch3.Merge() | Out-Null
$swatch3.Interior.Color = 0x5E5EDF

$desc3 = $ws.Range("F29:P30")
$desc3.Merge() | Out-Null
$desc3.Value = "Employee has no time-in and therefore, considered as absent."

Write-Host "edit complete"
